$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "91.616.24"
Set-TextValue "E2" "  +0.64%  "
Set-TextValue "D3" "3.116.85"
Set-TextValue "E3" "  +0.91%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "245.95"
Set-TextValue "E5" "  +0.54%  "
Set-TextValue "D6" "616.68"
Set-TextValue "E6" "  -0.39%  "
Set-TextValue "E7" "  -2.84%  "
Set-TextValue "D8" "0.383"
Set-TextValue "E8" "  +4.12%  "
Set-TextValue "E9" "  -0.08%  "
Set-TextValue "D10" "3.115.12"
Set-TextValue "E10" "  +0.98%  "
Set-TextValue "D11" "0.731"
Set-TextValue "E11" "  -2.20%  "
Set-TextValue "E12" "  +1.19%  "
Set-TextValue "D13" "0.0000250"
Set-TextValue "E13" "  +0.32%  "
Set-TextValue "B14" "Toncoin"
Set-TextValue "C14" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D14" "5.59"
Set-TextValue "E14" "  +2.11%  "
Set-TextValue "B15" "Avalanche"
Set-TextValue "C15" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D15" "34.73"
Set-TextValue "E15" "  -1.81%  "
Set-TextValue "D16" "91.563.96"
Set-TextValue "E16" "  +0.55%  "
Set-TextValue "D17" "3.706.45"
Set-TextValue "E17" "  +1.10%  "
Set-TextValue "D18" "3.081.61"
Set-TextValue "E18" "  -0.24%  "
Set-TextValue "D19" "3.70"
Set-TextValue "E19" "  +0.55%  "
Set-TextValue "D20" "14.81"
Set-TextValue "E20" "  +1.51%  "
Set-TextValue "D21" "5.80"
Set-TextValue "E21" "  -0.28%  "
Set-TextValue "D22" "9.47"
Set-TextValue "E22" "  +4.04%  "
Set-TextValue "D23" "446.42"
Set-TextValue "E23" "  +0.59%  "
Set-TextValue "D24" "0.0000201"
Set-TextValue "E24" "  -5.02%  "
Set-TextValue "D25" "5.84"
Set-TextValue "E25" "  +4.02%  "
Set-TextValue "D26" "87.87"
Set-TextValue "E26" "  -3.57%  "
Set-TextValue "D27" "11.72"
Set-TextValue "E27" "  -1.33%  "
Set-TextValue "D29" "0.144"
Set-TextValue "E29" "  +28.99%  "
Set-TextValue "D30" "1.00"
Set-TextValue "E30" "  +0.25%  "
Set-TextValue "D31" "0.234"
Set-TextValue "E31" "  -4.76%  "
Set-TextValue "E33" "  +4.23%  "
Set-TextValue "D34" "9.29"
Set-TextValue "E34" "  +0.38%  "
Set-TextValue "E35" "  -1.01%  "
Set-TextValue "D36" "7.84"
Set-TextValue "E36" "  +1.54%  "
Set-TextValue "D37" "26.20"
Set-TextValue "E37" "  -1.03%  "
Set-TextValue "D38" "4.19"
Set-TextValue "E38" "  +0.08%  "
Set-TextValue "D39" "1.94"
Set-TextValue "E39" "  +0.87%  "
Set-TextValue "D40" "490.22"
Set-TextValue "E40" "  -0.87%  "
Set-TextValue "D41" "1.30"
Set-TextValue "E41" "  +0.31%  "
Set-TextValue "D42" "0.438"
Set-TextValue "E42" "  +4.72%  "
Set-TextValue "D43" "3.40"
Set-TextValue "E43" "  -5.80%  "
Set-TextValue "D44" "22.18"
Set-TextValue "E44" "  +0.34%  "
Set-TextValue "D46" "157.83"
Set-TextValue "E46" "  +2.50%  "
Set-TextValue "D47" "0.705"
Set-TextValue "E47" "  +1.81%  "
Set-TextValue "D48" "1.90"
Set-TextValue "E48" "  +0.29%  "
Set-TextValue "E49" "  +1.29%  "
Set-TextValue "D50" "44.06"
Set-TextValue "E50" "  -0.23%  "
Set-TextValue "D51" "4.38"
Set-TextValue "E51" "  -2.25%  "
